$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 77, shifting existing rows 77..140 down to 78..141.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new record.
$ws.Cells.Item(77, 1).Value = 4
$ws.Cells.Item(77, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77, 3).Value = "Los Lagos"
$ws.Cells.Item(77, 4).Value = 44447
$ws.Cells.Item(77, 5).Value = 10
$ws.Cells.Item(77, 6).Value = 100112043
$ws.Cells.Item(77, 7).Value = "Pepino ensalada"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 80
$ws.Cells.Item(77, 11).Value = 21000
$ws.Cells.Item(77, 12).Value = 21000
$ws.Cells.Item(77, 13).Value = 21000
$ws.Cells.Item(77, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(77, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(77, 16).Value = 350
$ws.Cells.Item(77, 17).Value = 60
$ws.Cells.Item(77, 18).Value = "Hortaliza"
